# Re-process the metadata sheet with the newly curated dimensions.
# "municipio-nombre" (column C) moves from being a measure to a
# dimension (like provincia-nombre / comarca-nombre), while
# "gestion-explotacion" (column E) moves from being a dimension to a
# measure (like explotaciones) and therefore no longer needs a mapping
# file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (municipio-nombre): measure -> dimension
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# Column E (gestion-explotacion): dimension -> measure
$ws.Range("E2").Value = "iaest-measure:gestion-explotacion"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"

# The mapping file reference for gestion-explotacion is no longer
# needed now that it is a measure rather than a dimension.
$ws.Range("E5").Clear()
